# Mejoras en la parte de la interfaz
# Adds a new sheet "Hoja2" (positioned after "Hoja1") containing a 12x4
# grid of the letters a..l (same value repeated across columns A-D for
# each row), makes it the active sheet/tab, and leaves the selection on
# cell B14.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

$ws2 = $wb.Worksheets.Add($null, $sheet1)
$ws2.Name = "Hoja2"

$letters = @("a", "b", "c", "d", "e", "f", "g", "h", "i", "j", "k", "l")
# Row 6 ("f") is written before row 5 ("e") to match the original
# shared-string allocation order from the authored workbook.
$rowOrder = @(1, 2, 3, 4, 6, 5, 7, 8, 9, 10, 11, 12)

foreach ($row in $rowOrder) {
    $value = $letters[$row - 1]
    $ws2.Cells.Item($row, 1).Value = $value
    $ws2.Cells.Item($row, 2).Value = $value
    $ws2.Cells.Item($row, 3).Value = $value
    $ws2.Cells.Item($row, 4).Value = $value
}

$ws2.Range("B14").Select()
